$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.599.46"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "3.393.90"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.23"
$ws.Range("E5").Value = "  +1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.35"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.472"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.62"
$ws.Range("E9").Value = "  +0.24%  "
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.385"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("D12").Value = "3.970.05"
$ws.Range("E12").Value = "  +0.15%  "
$ws.Range("E13").Value = "  -0.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.00"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Value = "3.389.12"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "61.636.27"
$ws.Range("E17").Value = "  +0.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.13"
$ws.Range("E18").Value = "  +0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.68"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.13"
$ws.Range("E20").Value = "  +1.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "386.88"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "74.38"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.548"
$ws.Range("E23").Value = "  -0.54%  "
$ws.Range("E25").Value = "  -2.90%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.14"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.39"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "168.44"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  +2.67%  "
$ws.Range("D37").Value = "3.423.42"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "27.53"
$ws.Range("E39").Value = "  +2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0755"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.782"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.43"
$ws.Range("E42").Value = "  +1.10%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  +2.65%  "
$ws.Range("D45").Value = "2.475.23"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.72"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("E48").Value = "  -1.60%  "
$ws.Range("E49").Value = "  -0.40%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.03"
$ws.Range("E50").Value = "  -5.63%  "
$ws.Range("E51").Value = "  -1.24%  "
